$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BINARY")
Write-Host $ws.Name
